$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 30: fill in the remaining columns (B30 already has "4")
# Set C30's text first so it claims the shared-string slot formerly
# used by C29's old text (keeps the new shared-string order in line
# with the canonical save).
$ws.Range("C30").Value = "Format paths.json"
$ws.Range("D30").Value = 0.66666666666666663
$ws.Range("E30").Value = 0.75
$ws.Range("F30").Value = 2

# Row 29: update existing description text
$ws.Range("C29").Value = "Format nodes.json, Show nearest nodes inside circle."

# Row 31: brand new row
$ws.Range("B31").Value = "5"
$ws.Range("C31").Value = "Find nearest path to accident spot"
$ws.Range("D31").Value = 0.75
$ws.Range("E31").Value = 0.79166666666666663
$ws.Range("F31").Value = 1

$ws.Range("C32").Select()
